$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "OrderNo"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Fulfilment Store"
$ws.Range("D1").Value = "Total Price"

$ws.Range("D2").NumberFormat = """$""#,##0;[Red]\-""$""#,##0"

$ws.Range("A2:D2").Select()
